$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.415233185193614
$ws.Range("C2").Value = 0.2383185235071323
$ws.Range("D2").Value = 0.3914680602252503
$ws.Range("E2").Value = 0.08201193451269972
$ws.Range("G2").Value = 0.002730821100359145
$ws.Range("I2").Value = 3.760105428490434
$ws.Range("J2").Value = 0.01429883411345045
$ws.Range("K2").Value = 2.39086458928
$ws.Range("L2").Value = 0.6317356454288756
$ws.Range("M2").Value = 0.6434664213387506
$ws.Range("N2").Value = 5.925080115494893
$ws.Range("B3").Value = 2.389524415994202
$ws.Range("C3").Value = 0.2309697580486443
$ws.Range("D3").Value = 0.3907700879249347
$ws.Range("E3").Value = 0.08227695199060303
$ws.Range("G3").Value = 0.002735772874106025
$ws.Range("I3").Value = 3.751328576363932
$ws.Range("J3").Value = 0.01358858181872691
$ws.Range("K3").Value = 2.35787048725075
$ws.Range("L3").Value = 0.6305845061558841
$ws.Range("M3").Value = 0.6390979676655562
$ws.Range("N3").Value = 5.919947293002451
$ws.Range("B4").Value = 2.375007363536241
$ws.Range("C4").Value = 0.2265926391253004
$ws.Range("D4").Value = 0.3904934815317631
$ws.Range("E4").Value = 0.08245716340387599
$ws.Range("G4").Value = 0.002738975076443583
$ws.Range("I4").Value = 3.746951151975409
$ws.Range("J4").Value = 0.01314645017480487
$ws.Range("K4").Value = 2.338898486097435
$ws.Range("L4").Value = 0.6301170709420347
$ws.Range("M4").Value = 0.6367199781001744
$ws.Range("N4").Value = 5.917686977726007
$ws.Range("B5").Value = 2.369410588133206
$ws.Range("C5").Value = 0.2248427713347496
$ws.Range("D5").Value = 0.3904190194196531
$ws.Range("E5").Value = 0.08253501030780441
$ws.Range("G5").Value = 0.002740320817959985
$ws.Range("I5").Value = 3.745421742804027
$ws.Range("J5").Value = 0.01296474381947377
$ws.Range("K5").Value = 2.331490688639576
$ws.Range("L5").Value = 0.6299868384497316
$ws.Range("M5").Value = 0.6358274853006307
$ws.Range("N5").Value = 5.916989906162982
$ws.Range("B6").Value = 2.368500521344401
$ws.Range("C6").Value = 0.2245542478372045
$ws.Range("D6").Value = 0.3904089671143254
$ws.Range("E6").Value = 0.08254820337305802
$ws.Range("G6").Value = 0.002740546746430583
$ws.Range("I6").Value = 3.745183154022214
$ws.Range("J6").Value = 0.01293447846822282
$ws.Range("K6").Value = 2.330280163583836
$ws.Range("L6").Value = 0.6299688542626072
$ws.Range("M6").Value = 0.6356839129189638
$ws.Range("N6").Value = 5.916887686041463
$ws.Range("B7").Value = 2.374930591503784
$ws.Range("C7").Value = 0.2265689029168243
$ws.Range("D7").Value = 0.3904923223431354
$ws.Range("E7").Value = 0.08245819540594024
$ws.Range("G7").Value = 0.00273899306034462
$ws.Range("I7").Value = 3.74692949561188
$ws.Range("J7").Value = 0.01314400585344444
$ws.Range("K7").Value = 2.338797272398551
$ws.Range("L7").Value = 0.6301150705458411
$ws.Range("M7").Value = 0.6367076315875266
$ws.Range("N7").Value = 5.917676669857769
$ws.Range("B8").Value = 2.406105641426024
$ws.Range("C8").Value = 0.2357565656562315
$ws.Range("D8").Value = 0.3911958892372667
$ws.Range("E8").Value = 0.08209968936011514
$ws.Range("G8").Value = 0.002732494972490339
$ws.Range("I8").Value = 3.756869233987899
$ws.Range("J8").Value = 0.01405518114083293
$ws.Range("K8").Value = 2.379221153337795
$ws.Range("L8").Value = 0.6312890854045037
$ws.Range("M8").Value = 0.6418970553495953
$ws.Range("N8").Value = 5.923125276453504
$ws.Range("B9").Value = 2.477303745054883
$ws.Range("C9").Value = 0.254851637028878
$ws.Range("D9").Value = 0.3937796472056476
$ws.Range("E9").Value = 0.08153494325874444
$ws.Range("G9").Value = 0.002721029903768102
$ws.Range("I9").Value = 3.784389261517589
$ws.Range("J9").Value = 0.01579494035311768
$ws.Range("K9").Value = 2.468710606763352
$ws.Range("L9").Value = 0.6354887834626197
$ws.Range("M9").Value = 0.6544869171965857
$ws.Range("N9").Value = 5.940888784238467
$ws.Range("B10").Value = 2.535759877843702
$ws.Range("C10").Value = 0.2695488496148357
$ws.Range("D10").Value = 0.3964105550820562
$ws.Range("E10").Value = 0.08120368438787473
$ws.Range("G10").Value = 0.002713376875372282
$ws.Range("I10").Value = 3.809510541820288
$ws.Range("J10").Value = 0.01704568686866637
$ws.Range("K10").Value = 2.540713803487506
$ws.Range("L10").Value = 0.6397297826324433
$ws.Range("M10").Value = 0.6652090788465372
$ws.Range("N10").Value = 5.958269283962863
$ws.Range("B11").Value = 2.563691174437679
$ws.Range("C11").Value = 0.2763825277848184
$ws.Range("D11").Value = 0.3977662738965932
$ws.Range("E11").Value = 0.08107101904328218
$ws.Range("G11").Value = 0.002710060769834278
$ws.Range("I11").Value = 3.822005596589548
$ws.Range("J11").Value = 0.01760899248397152
$ws.Range("K11").Value = 2.574834707441454
$ws.Range("L11").Value = 0.6419098685100408
$ws.Range("M11").Value = 0.6704070022963577
$ws.Range("N11").Value = 5.967119712353536
$ws.Range("B12").Value = 2.574460718210275
$ws.Range("C12").Value = 0.2789916896160776
$ws.Range("D12").Value = 0.3983024690668202
$ws.Range("E12").Value = 0.08102336337637794
$ws.Range("G12").Value = 0.002708828677751022
$ws.Range("I12").Value = 3.826890709020006
$ws.Range("J12").Value = 0.01782150672204708
$ws.Range("K12").Value = 2.587952203701377
$ws.Range("L12").Value = 0.642771448624444
$ws.Range("M12").Value = 0.6724213869030109
$ws.Range("N12").Value = 5.97060708910476
$ws.Range("B13").Value = 2.572132739988717
$ws.Range("C13").Value = 0.2784288061289999
$ws.Range("D13").Value = 0.3981859757603416
$ws.Range("E13").Value = 0.08103351221268795
$ws.Range("G13").Value = 0.002709092981422591
$ws.Range("I13").Value = 3.825831784755351
$ws.Range("J13").Value = 0.01777577316613943
$ws.Range("K13").Value = 2.585118367087034
$ws.Range("L13").Value = 0.6425842901644927
$ws.Range("M13").Value = 0.6719855056674646
$ws.Range("N13").Value = 5.969849972597558
$ws.Range("B14").Value = 2.564573332025532
$ws.Range("C14").Value = 0.2765967554959161
$ws.Range("D14").Value = 0.3978099299765745
$ws.Range("E14").Value = 0.08106704669808273
$ws.Range("G14").Value = 0.00270995893167203
$ws.Range("I14").Value = 3.822404421001082
$ws.Range("J14").Value = 0.0176264920398701
$ws.Range("K14").Value = 2.575909949507178
$ws.Range("L14").Value = 0.6419800293783595
$ws.Range("M14").Value = 0.6705718042678157
$ws.Range("N14").Value = 5.967403896059352
$ws.Range("B15").Value = 2.559968053931073
$ws.Range("C15").Value = 0.2754773620964954
$ws.Range("D15").Value = 0.3975825610513368
$ws.Range("E15").Value = 0.08108792346405025
$ws.Range("G15").Value = 0.002710492427164075
$ws.Range("I15").Value = 3.820325053993216
$ws.Range("J15").Value = 0.0175349497117665
$ws.Range("K15").Value = 2.570295143114208
$ws.Range("L15").Value = 0.6416145930495674
$ws.Range("M15").Value = 0.669711866957492
$ws.Range("N15").Value = 5.965923308470224
$ws.Range("B16").Value = 2.533961449988794
$ws.Range("C16").Value = 0.2691052392745235
$ws.Range("D16").Value = 0.3963251501069891
$ws.Range("E16").Value = 0.08121271622646464
$ws.Range("G16").Value = 0.002713596906063875
$ws.Range("I16").Value = 3.808715439465701
$ws.Range("J16").Value = 0.01700876099864246
$ws.Range("K16").Value = 2.538511430449915
$ws.Range("L16").Value = 0.6395923527918086
$ws.Range("M16").Value = 0.6648758266966297
$ws.Range("N16").Value = 5.957709896980077
$ws.Range("B17").Value = 2.518350249613206
$ws.Range("C17").Value = 0.2652341108531289
$ws.Range("D17").Value = 0.3955944423976518
$ws.Range("E17").Value = 0.0812938820721012
$ws.Range("G17").Value = 0.002715543650929849
$ws.Range("I17").Value = 3.801866690553055
$ws.Range("J17").Value = 0.01668452253735708
$ws.Range("K17").Value = 2.519363215123747
$ws.Range("L17").Value = 0.6384159889445016
$ws.Range("M17").Value = 0.6619911066221249
$ws.Range("N17").Value = 5.952913124534774
$ws.Range("B18").Value = 2.50949716276449
$ws.Range("C18").Value = 0.2630214475023536
$ws.Range("D18").Value = 0.3951891177175071
$ws.Range("E18").Value = 0.08134226381916143
$ws.Range("G18").Value = 0.002716678932909354
$ws.Range("I18").Value = 3.798027924134615
$ws.Range("J18").Value = 0.01649749487268437
$ws.Range("K18").Value = 2.508478258108426
$ws.Range("L18").Value = 0.6377629875881752
$ws.Range("M18").Value = 0.6603620468086859
$ws.Range("N18").Value = 5.950242986456658
$ws.Range("B19").Value = 2.506521311954941
$ws.Range("C19").Value = 0.2622746621789247
$ws.Range("D19").Value = 0.3950544521106281
$ws.Range("E19").Value = 0.08135893687563378
$ws.Range("G19").Value = 0.002717065996991306
$ws.Range("I19").Value = 3.796745435667361
$ws.Range("J19").Value = 0.01643407823066667
$ws.Range("K19").Value = 2.504814878878676
$ws.Range("L19").Value = 0.6375459492050197
$ws.Range("M19").Value = 0.6598156550287158
$ws.Range("N19").Value = 5.949354176515129
$ws.Range("B20").Value = 2.519999041763015
$ws.Range("C20").Value = 0.2656447587438322
$ws.Range("D20").Value = 0.3956706795415386
$ws.Range("E20").Value = 0.08128506624005372
$ws.Range("G20").Value = 0.002715334806428887
$ws.Range("I20").Value = 3.802585354889956
$ws.Range("J20").Value = 0.01671909344558031
$ws.Range("K20").Value = 2.521388265521296
$ws.Range("L20").Value = 0.6385387714339856
$ws.Range("M20").Value = 0.6622950692664844
$ws.Range("N20").Value = 5.953414553533008
$ws.Range("B21").Value = 2.56678848727006
$ws.Range("C21").Value = 0.2771342913894159
$ws.Range("D21").Value = 0.3979197648627917
$ws.Range("E21").Value = 0.08105712681850008
$ws.Range("G21").Value = 0.00270970394033013
$ws.Range("I21").Value = 3.823406954371762
$ws.Range("J21").Value = 0.01767036102743447
$ws.Range("K21").Value = 2.578609346427726
$ws.Range("L21").Value = 0.6421565379197034
$ws.Range("M21").Value = 0.6709857931941556
$ws.Range("N21").Value = 5.968118678023586
$ws.Range("B22").Value = 2.598490526188186
$ws.Range("C22").Value = 0.2847681436968799
$ws.Range("D22").Value = 0.3995226327186288
$ws.Range("E22").Value = 0.08092319934823422
$ws.Range("G22").Value = 0.002706161607910506
$ws.Range("I22").Value = 3.837909848070723
$ws.Range("J22").Value = 0.0182874297526574
$ws.Range("K22").Value = 2.617152953199422
$ws.Range("L22").Value = 0.6447309439362527
$ws.Range("M22").Value = 0.6769340533429684
$ws.Range("N22").Value = 5.978520954636082
$ws.Range("B23").Value = 2.581467853199854
$ws.Range("C23").Value = 0.2806823500766598
$ws.Range("D23").Value = 0.3986549978602199
$ws.Range("E23").Value = 0.08099330579569042
$ws.Range("G23").Value = 0.00270803965289192
$ws.Range("I23").Value = 3.830087494907644
$ws.Range("J23").Value = 0.01795850735913973
$ws.Range("K23").Value = 2.596476553292518
$ws.Range("L23").Value = 0.6433377341817561
$ws.Range("M23").Value = 0.673734806104143
$ws.Range("N23").Value = 5.972896511444389
$ws.Range("B24").Value = 2.519253242814329
$ws.Range("C24").Value = 0.2654590646626218
$ws.Range("D24").Value = 0.3956361667238326
$ws.Range("E24").Value = 0.08128904652441271
$ws.Range("G24").Value = 0.002715429175024865
$ws.Range("I24").Value = 3.802260139386746
$ws.Range("J24").Value = 0.01670346586614357
$ws.Range("K24").Value = 2.520472355230567
$ws.Range("L24").Value = 0.6384831888735079
$ws.Range("M24").Value = 0.6621575561632085
$ws.Range("N24").Value = 5.953187584637249
$ws.Range("B25").Value = 2.456964199429564
$ws.Range("C25").Value = 0.2495693852971215
$ws.Range("D25").Value = 0.3929519185076344
$ws.Range("E25").Value = 0.0816729845277111
$ws.Range("G25").Value = 0.002723995611647085
$ws.Range("I25").Value = 3.776084078833549
$ws.Range("J25").Value = 0.01532919977552538
$ws.Range("K25").Value = 2.443404763897036
$ws.Range("L25").Value = 0.6341496549592307
$ws.Range("M25").Value = 0.6508225485438501
$ws.Range("N25").Value = 5.935324001718868
